# Inclusion/Exclusion template update:
# - rename Sheet1 -> exclusion_rules
# - add new sheet "store_policy"
# - add new exclusion rule rows (Hero SKU / Brand / Sub Brand / PepsiCo Segment / PepsiCo Sub Segment
#   Space to Sales Index KPIs excluded by brand_name) to exclusion_rules
# - populate store_policy with the additional_attribute_1 / store_type mapping for the same KPIs
# - tweak the "additional display,stock" value down to just "additional display"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "exclusion_rules"

# --- fix up existing value on exclusion_rules (row4 / D4: "additional display,stock" -> "additional display")
$ws1.Range("D4").Value2 = "additional display"

$brandList = 'BLACK COUNTRY SNACKS, AKSAM PALUSZKI, LAJKONIK PALUSZKI, ALKA ELEPHANT, RUMPLERS, TABITHA, CHEF''S LARDER, BOROMIR, COFRESH, JIFFY POP, CROCO, THE CURATORS, FRESHERS, FUDCO, GEFEN, GINNI''S, OH MY GURU!, HALDIRAMS, HALDIRAMS SNACKS, INDIE BAY SNACKS, INNATE, JACK-LNK''S, EAZY-PP-PPCRN, ZWEIFEL CRISPS, HUNKY DORYS CRISPS, LAJKONIK JUNIOR, LOVE CHIN CHIN, NISHAS SNACKS, NUTELLA, OSEM SAVOURY SNACK, OUR LITTLE REBELL!ON, EPIC, CRAWFORDS, FLIPZ, OATIS, RYMUT SNACKS, GINCO, SUNSHINE SNACKS, JAY''S, MIDLAND SNACKS, RED MILL SAVOURY SNACKS, SENSIBLE PORTIONS, VISCONTI SNACKS, WELL & TRULY SNACKS, WILD WEST, WILDING''S, BLUE DRAGON, BEPPS, BLUE DIAMOND, COFRESH SNACKS, SCHAR, OLD EL PASO, PLANTERS, LINWOODS, CYPRESSA, KOHINOOR SNACKS, KOIKEYA, PALUSZKI, LORENZ CRISPS, MCCOLGAN, ITSU, NAIRNS, NATURES STORE SNACKS, NIM''S, BAMBA SNACKS, BISSLI SNACKS, SHARWOODS, MR PORKY SNACKS, TYGRYSKI, THE REAL PORK CRACKLING CO SNACKS, THE SNAFFLING PIG CO, WHITWORTHS, YUM & YAY'

$kpis = @(
    "Hero SKU Space to Sales Index",
    "Brand Space to Sales Index",
    "Sub Brand Space to Sales Index",
    "PepsiCo Segment Space to Sales Index",
    "PepsiCo Sub Segment Space to Sales Index"
)

# --- new exclusion rule rows on exclusion_rules (rows 5-9): Exclude brand_name = <brandList>
$row = 5
foreach ($kpi in $kpis) {
    $ws1.Range("A" + $row).Value2 = $kpi
    $ws1.Range("B" + $row).Value2 = "Exclude"
    $ws1.Range("C" + $row).Value2 = "brand_name"
    $ws1.Range("D" + $row).Value2 = $brandList
    $row = $row + 1
}

# leave row 10 untouched/blank, row 11 stays blank as well (dimension extends to row 11)

# resize columns on exclusion_rules to fit the new wide content
$ws1.Columns.Item(1).ColumnWidth = 39.3117408906883
$ws1.Columns.Item(2).ColumnWidth = 7.71255060728745
$ws1.Columns.Item(3).ColumnWidth = 21.1012145748988
$ws1.Columns.Item(4).ColumnWidth = 51.2024291497976
$ws1.Columns.Item(5).ColumnWidth = 22.1740890688259
$ws1.Columns.Item(6).ColumnWidth = 21.1012145748988

$ws1.Range("G5").Select()

# --- add the store_policy sheet right after exclusion_rules
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "store_policy"

$ws2.Range("A1").Value2 = "KPI"
$ws2.Range("B1").Value2 = "additional_attribute_1"
$ws2.Range("C1").Value2 = "store_type"

# match the bold/grey header style already used for exclusion_rules!A1
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)

# match the blue header fill (fontId 0 / fillId matching the existing blue fill) via a
# scratch cell + format-only paste, so no stray intermediate style is left behind
$ws2.Range("Z1").Interior.ColorIndex = 37
$ws2.Range("Z1").Interior.PatternColorIndex = 15
$ws2.Range("Z1").Copy()
$ws2.Range("B1:C1").PasteSpecial(-4122)
$ws2.Range("Z1").Clear()

$row = 2
foreach ($kpi in $kpis) {
    $ws2.Range("A" + $row).Value2 = $kpi
    $ws2.Range("B" + $row).Value2 = "TT"
    $row = $row + 1
}

$ws2.Columns.Item(1).ColumnWidth = 38.8825910931174
$ws2.Columns.Item(2).ColumnWidth = 21.1012145748988
$ws2.Columns.Item(3).ColumnWidth = 13.0688259109312

$ws2.Range("E9").Select()

# activate the exclusion_rules sheet (first tab) as the active tab, matching the source
$ws1.Activate()
$ws1.Range("G5").Select()
